# Reconcile treatment support code with new approach to treatment outcomes
#
# The "ideal" treatment-death parameter (int_prop_treatment_death_ideal) is
# no longer used now that treatment outcomes are derived differently, so its
# row is removed from the "constants" sheet. All rows below it shift up by
# one; Excel re-points every downstream formula/shared-string reference
# automatically when we delete the whole row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Row 58 holds: A58 = "int_prop_treatment_death_ideal", B58 = 0.016
$ws.Rows.Item(58).Delete()

# Restore the user's on-screen selection to reflect the new layout
# (previously selected row 62 -> now row 55 after the row shift context).
$ws.Activate()
$ws.Range("A55").Select()
